# Insert a new weekly price record for "Locoto" (Vega Modelo de Temuco)
# before the existing row 13, shifting the remaining historical rows
# down by one (13->14, 14->15, ..., 29->30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 13 (pushes old rows 13..29 to 14..30)
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the latest weekly entry
$ws.Cells.Item(13, 1).Value = 10
$ws.Cells.Item(13, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(13, 3).Value = "La Araucanía"
$ws.Cells.Item(13, 4).Value = 44757
$ws.Cells.Item(13, 5).Value = 9
$ws.Cells.Item(13, 6).Value = 100112042
$ws.Cells.Item(13, 7).Value = "Locoto"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 80
$ws.Cells.Item(13, 11).Value = 3300
$ws.Cells.Item(13, 12).Value = 3300
$ws.Cells.Item(13, 13).Value = 3300
$ws.Cells.Item(13, 14).Value = "`$/kilo"
$ws.Cells.Item(13, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(13, 16).Value = 3300
$ws.Cells.Item(13, 17).Value = 1
$ws.Cells.Item(13, 18).Value = "Hortaliza"

# Match the date number format used by the other rows in column D
$ws.Cells.Item(13, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
